$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) storage for Price/Volume columns so values like
# "1.00", "34.00", "0.0000262", "88.925.21" are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '88.925.21'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '3.272.12'
$ws.Range("E3").Value = '  -2.50%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '212.12'
$ws.Range("E5").Value = '  -3.73%  '
$ws.Range("D6").Value = '625.48'
$ws.Range("E6").Value = '  -2.58%  '
$ws.Range("D7").Value = '0.374'
$ws.Range("E7").Value = '  +14.02%  '
$ws.Range("D8").Value = '0.715'
$ws.Range("E8").Value = '  +15.54%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '3.266.67'
$ws.Range("E10").Value = '  -2.90%  '
$ws.Range("D11").Value = '0.574'
$ws.Range("E11").Value = '  -6.25%  '
$ws.Range("E12").Value = '  +11.13%  '
$ws.Range("D13").Value = '0.0000262'
$ws.Range("E13").Value = '  -5.80%  '
$ws.Range("D14").Value = '5.46'
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '3.873.89'
$ws.Range("E15").Value = '  -2.48%  '
$ws.Range("D16").Value = '34.00'
$ws.Range("E16").Value = '  -1.46%  '
$ws.Range("D17").Value = '88.798.44'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '3.294.02'
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("D19").Value = '3.18'
$ws.Range("E19").Value = '  -1.54%  '
$ws.Range("D20").Value = '14.05'
$ws.Range("E20").Value = '  -4.71%  '
$ws.Range("D21").Value = '436.37'
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").Value = '8.87'
$ws.Range("E22").Value = '  -3.18%  '
$ws.Range("D23").Value = '5.33'
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").Value = '7.38'
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("D25").Value = '5.22'
$ws.Range("E25").Value = '  -4.22%  '
$ws.Range("D26").Value = '12.22'
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("D27").Value = '3.461.89'
$ws.Range("E27").Value = '  -2.01%  '
$ws.Range("D28").Value = '76.90'
$ws.Range("E28").Value = '  -2.71%  '
$ws.Range("E29").Value = '  +2.73%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("D31").Value = '0.180'
$ws.Range("E31").Value = '  -5.12%  '
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").Value = '562.54'
$ws.Range("E33").Value = '  -7.06%  '
$ws.Range("D34").Value = '8.72'
$ws.Range("E34").Value = '  -6.62%  '
$ws.Range("D35").Value = '1.38'
$ws.Range("E35").Value = '  -11.91%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '7.13'
$ws.Range("E36").Value = '  +5.96%  '
$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").Value = '1.96'
$ws.Range("E37").Value = '  -5.19%  '
$ws.Range("E38").Value = '  -8.45%  '
$ws.Range("D39").Value = '22.66'
$ws.Range("E39").Value = '  -3.76%  '
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").Value = '3.09'
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D43").Value = '0.401'
$ws.Range("E43").Value = '  -4.95%  '
$ws.Range("D44").Value = '2.02'
$ws.Range("E44").Value = '  -2.24%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").Value = '153.59'
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("D47").Value = '180.52'
$ws.Range("E47").Value = '  -5.14%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.133'
$ws.Range("E48").Value = '  +17.67%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '44.71'
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("E50").Value = '  -5.29%  '
$ws.Range("D51").Value = '4.21'
$ws.Range("E51").Value = '  -2.46%  '

# Restore the default (General) style on the edited cells so no stray
# number-format styling is left behind, while keeping the values as text.
$ws.Range("D2:E51").Style = "Normal"
